# Update the date line at the top of the document.
$d = $word.ActiveDocument
$null = $d.Content.Find.Execute("2025-01-22 Wednesday", $true, $false, $false, $false, $false,
                                 $true, 1, $false, "2025-01-23 Thursday", 2)

# The worksheet table: data lives in table rows 1, 5, 9, 13, 17 (each
# followed by 3 blank spacer rows), 5 columns per row. Address every
# cell directly by (row, column) so duplicate text elsewhere in the
# table can't cause an ambiguous Find/Replace match.
$t = $d.Tables.Item(1)

$grid = @{
    1  = @("23÷6=3, 5",  "23÷5=4, 3",  "63÷5=12, 3", "24÷5=4, 4",  "53÷5=10, 3")
    5  = @("32÷4=8, 0",  "73÷5=14, 3", "62÷7=8, 6",  "73÷2=36, 1", "91÷5=18, 1")
    9  = @("93÷5=18, 3", "18÷5=3, 3",  "29÷4=7, 1",  "58÷7=8, 2",  "37÷3=12, 1")
    13 = @("36÷9=4, 0",  "73÷3=24, 1", "96÷7=13, 5", "23÷2=11, 1", "43÷9=4, 7")
    17 = @("31÷3=10, 1", "92÷6=15, 2", "74÷2=37, 0", "62÷8=7, 6",  "59÷2=29, 1")
}

foreach ($rowIndex in $grid.Keys) {
    $values = $grid[$rowIndex]
    for ($col = 1; $col -le $values.Count; $col++) {
        $cell = $t.Cell($rowIndex, $col)
        $cell.Range.Text = $values[$col - 1]
    }
}
